$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 ("Live Châtel" -> "Châtel SKI resort", Châtel / France) gains the
# CITY / COUNTRY values that every other row already has, and its Location
# name is reworded.
$ws.Range("D12").Value = "Châtel SKI resort"
$ws.Range("E12").Value = "Châtel"
$ws.Range("F12").Value = "France"

# Move the active selection, matching the saved view state in the workbook.
$ws.Range("F16").Select()
